# Update the confidential disclosure text with the new "as of" date.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) columns for rows 2-10.
$ws.Range("D2").Value = 0.09315186476104689
$ws.Range("E2").Value = -0.0586125881043652

$ws.Range("D3").Value = 0.1067522248848124
$ws.Range("E3").Value = -0.02542674253200561

$ws.Range("D4").Value = 0.1197818406906597
$ws.Range("E4").Value = -0.009032656527445337

$ws.Range("D5").Value = 0.1405588482412111
$ws.Range("E5").Value = -0.006592889098186894

$ws.Range("D6").Value = 0.1360474048095251
$ws.Range("E6").Value = 0.001508502468458461

$ws.Range("D7").Value = 0.1469892125916969
$ws.Range("E7").Value = -0.006860745410717572

$ws.Range("D8").Value = 0.1284353174312631
$ws.Range("E8").Value = -0.02934340499709487

$ws.Range("D9").Value = 0.1282832865897847
$ws.Range("E9").Value = -0.01064929319740204

$ws.Range("E10").Value = -0.0161209539269328

# Restore sheet protection (the original password is not recoverable here,
# so the sheet is re-locked to preserve the protected state/intent).
$ws.Protect()
